# Update Leave Card 4/27/2023 4:55 PM
#
# A new Leave Card entry pair is inserted into the "Table1" table on Sheet1:
#   - Existing row 163 (PERIOD 4/1/2023) gets a "SP(1-0-0)" particular with a
#     REMARKS date of 4/25/2023.
#   - A brand-new row is inserted right after it (becomes the new row 164)
#     recording a "SL(1-0-0)" particular, 1 day used, with a REMARKS date of
#     4/13/2023.
# Inserting that row inside the table shifts every following data row (and
# the trailing "last row" of the table) down by one, and the table grows
# from A8:K194 to A8:K195.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$lo = $ws.ListObjects.Item("Table1")

# Insert a new worksheet row at row 164 (this is the row right after the
# period-starting row 163). Everything below shifts down by one row.
$ws.Rows("164:164").Insert()

# The table's own range does not auto-grow when a row is inserted via
# Rows.Insert, so extend it explicitly to include the new last row (195).
$lo.Resize($ws.Range("A8:K195"))

# The newly inserted worksheet row (164) comes back with "no border" default
# formatting instead of the table's normal row look. Copy the formatting
# (borders, number formats, fonts, etc.) from the row directly below it
# (165), which still holds the formatting that used to belong to the old
# row 164 before the insert shifted it down.
$ws.Range("A165:K165").Copy()
$ws.Range("A164:K164").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Restore the calculated "EARNED " helper-column formula for the new row,
# matching the rest of the table.
$ws.Range("G164").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
# The row that now sits at the very bottom of the table (195, formerly 194)
# also needs that helper formula restored/refreshed.
$ws.Range("G195").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Give K163/K164 (REMARKS column) the same date-formatted look already used
# elsewhere in the table (e.g. K162), then fill in the data for both rows.
$ws.Range("K162").Copy()
$ws.Range("K163").PasteSpecial(-4122)        # xlPasteFormats
$ws.Range("K164").PasteSpecial(-4122)        # xlPasteFormats
$excel.CutCopyMode = $false

# Row 163 (PERIOD 4/1/2023): particular "SP(1-0-0)", remarks date 4/25/2023.
$ws.Range("B163").Value = "SP(1-0-0)"
$ws.Range("K163").Value2 = 45041

# Row 164 (new row, no PERIOD date): particular "SL(1-0-0)", 1 day used,
# remarks date 4/13/2023.
$ws.Range("B164").Value = "SL(1-0-0)"
$ws.Range("H164").Value = 1
$ws.Range("K164").Value2 = 45029

# Match the author's final active-cell selection.
[void]$ws.Range("B165").Select()
